$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.45077905506912
$ws.Range("C2").Value = 0.2032878802410494
$ws.Range("D2").Value = 0.4437987674050845
$ws.Range("E2").Value = 0.1198194495716205
$ws.Range("G2").Value = 2.560828586082749
$ws.Range("H2").Value = 2.043167260661107
$ws.Range("I2").Value = 1.890953932940029
$ws.Range("J2").Value = 0.03463636307316342
$ws.Range("L2").Value = 0.6288189362564225
$ws.Range("B3").Value = 2.345194687835829
$ws.Range("C3").Value = 0.1806486192394061
$ws.Range("D3").Value = 0.4431670987413128
$ws.Range("E3").Value = 0.1203823107775754
$ws.Range("G3").Value = 2.551899791039659
$ws.Range("H3").Value = 2.046860853131591
$ws.Range("I3").Value = 1.897701294535729
$ws.Range("J3").Value = 0.03423350640571243
$ws.Range("L3").Value = 0.6206341304872467
$ws.Range("B4").Value = 2.28150492424561
$ws.Range("C4").Value = 0.1667957781242251
$ws.Range("D4").Value = 0.4429665283869042
$ws.Range("E4").Value = 0.1207619271742351
$ws.Range("G4").Value = 2.548007462160143
$ws.Range("H4").Value = 2.050160513903705
$ws.Range("I4").Value = 1.902862952720348
$ws.Range("J4").Value = 0.03398184900739132
$ws.Range("L4").Value = 0.6159107783568487
$ws.Range("B5").Value = 2.255837911607273
$ws.Range("C5").Value = 0.1611625147778
$ws.Range("D5").Value = 0.4429318942588054
$ws.Range("E5").Value = 0.1209251886249083
$ws.Range("G5").Value = 2.546819693806498
$ws.Range("H5").Value = 2.051764110760217
$ws.Range("I5").Value = 1.905221886798564
$ws.Range("J5").Value = 0.03387821051555306
$ws.Range("L5").Value = 0.6140619494108535
$ws.Range("B6").Value = 2.251593280070551
$ws.Range("C6").Value = 0.1602278297683313
$ws.Range("D6").Value = 0.4429289882687328
$ws.Range("E6").Value = 0.1209528156684829
$ws.Range("G6").Value = 2.546646485765052
$ws.Range("H6").Value = 2.052046012962421
$ws.Range("I6").Value = 1.905628999412123
$ws.Range("J6").Value = 0.03386093570253301
$ws.Range("L6").Value = 0.6137595415167851
$ws.Range("B7").Value = 2.281157606795034
$ws.Range("C7").Value = 0.1667197580154038
$ws.Range("D7").Value = 0.442965870581375
$ws.Range("E7").Value = 0.1207640942828014
$ws.Range("G7").Value = 2.547989832349913
$ws.Range("H7").Value = 2.050181092728479
$ws.Range("I7").Value = 1.90289373243931
$ws.Range("J7").Value = 0.03398045570535047
$ws.Range("L7").Value = 0.6158855367816614
$ws.Range("B8").Value = 2.414137352505804
$ws.Range("C8").Value = 0.195471865945251
$ws.Range("D8").Value = 0.4435421017395953
$ws.Range("E8").Value = 0.1200064702160724
$ws.Range("G8").Value = 2.557419065084076
$ws.Range("H8").Value = 2.044226364750614
$ws.Range("I8").Value = 1.893068582902892
$ws.Range("J8").Value = 0.03449834819823394
$ws.Range("L8").Value = 0.6259340963883773
$ws.Range("B9").Value = 2.68394494584777
$ws.Range("C9").Value = 0.2522425119192917
$ws.Range("D9").Value = 0.4461583286881847
$ws.Range("E9").Value = 0.1187902291134009
$ws.Range("G9").Value = 2.588594437314896
$ws.Range("H9").Value = 2.040760030114683
$ws.Range("I9").Value = 1.881915832347417
$ws.Range("J9").Value = 0.03548003766637819
$ws.Range("L9").Value = 0.648038876350796
$ws.Range("B10").Value = 2.887698414820591
$ws.Range("C10").Value = 0.2942053157436817
$ws.Range("D10").Value = 0.4489878244571628
$ws.Range("E10").Value = 0.1180603397549174
$ws.Range("G10").Value = 2.619333625469096
$ws.Range("H10").Value = 2.043255237070866
$ws.Range("I10").Value = 1.878713569733065
$ws.Range("J10").Value = 0.03618100688560943
$ws.Range("L10").Value = 0.6657479873469185
$ws.Range("B11").Value = 2.981597243645695
$ws.Range("C11").Value = 0.3133543424183074
$ws.Range("D11").Value = 0.4504724267976599
$ws.Range("E11").Value = 0.1177637172498489
$ws.Range("G11").Value = 2.63504137815184
$ws.Range("H11").Value = 2.045493309029268
$ws.Range("I11").Value = 1.878350570141436
$ws.Range("J11").Value = 0.03649558039734657
$ws.Range("L11").Value = 0.674124767461052
$ws.Range("B12").Value = 3.017328377301396
$ws.Range("C12").Value = 0.3206144416898837
$ws.Range("D12").Value = 0.4510630176499575
$ws.Range("E12").Value = 0.1176564761608176
$ws.Range("G12").Value = 2.641239180195896
$ws.Range("H12").Value = 2.046500044412255
$ws.Range("I12").Value = 1.878371173210212
$ws.Range("J12").Value = 0.03661408864945059
$ws.Range("L12").Value = 0.6773430465172794
$ws.Range("B13").Value = 3.009625318965902
$ws.Range("C13").Value = 0.3190504563678189
$ws.Range("D13").Value = 0.4509345598659991
$ws.Range("E13").Value = 0.1176793464904016
$ws.Range("G13").Value = 2.639893244282746
$ws.Range("H13").Value = 2.046276134389899
$ws.Range("I13").Value = 1.878359693469747
$ws.Range("J13").Value = 0.03658859300694672
$ws.Range("L13").Value = 0.6766478782189864
$ws.Range("B14").Value = 2.98453338652871
$ws.Range("C14").Value = 0.3139514582803713
$ws.Range("D14").Value = 0.4505204457315841
$ws.Range("E14").Value = 0.1177547926278617
$ws.Range("G14").Value = 2.635546262809527
$ws.Range("H14").Value = 2.045572939063192
$ws.Range("I14").Value = 1.878349092819477
$ws.Range("J14").Value = 0.03650534241988979
$ws.Range("L14").Value = 0.6743886115613407
$ws.Range("B15").Value = 2.969186462309835
$ws.Range("C15").Value = 0.3108293227999468
$ws.Range("D15").Value = 0.4502704882848576
$ws.Range("E15").Value = 0.1178016673363178
$ws.Range("G15").Value = 2.632916168108636
$ws.Range("H15").Value = 2.04516296601571
$ws.Range("I15").Value = 1.878363207557001
$ws.Range("J15").Value = 0.03645426921896089
$ws.Range("L15").Value = 0.6730107600551065
$ws.Range("B16").Value = 2.881586228966853
$ws.Range("C16").Value = 0.2929551013737068
$ws.Range("D16").Value = 0.4488947763144608
$ws.Range("E16").Value = 0.1180804365256325
$ws.Range("G16").Value = 2.618341926535237
$ws.Range("H16").Value = 2.043131221390411
$ws.Range("I16").Value = 1.878759368271801
$ws.Range("J16").Value = 0.03616036260835642
$ws.Range("L16").Value = 0.6652070041220526
$ws.Range("B17").Value = 2.828156052146255
$ws.Range("C17").Value = 0.2820053079680349
$ws.Range("D17").Value = 0.4481014006826598
$ws.Range("E17").Value = 0.1182605153101921
$ws.Range("G17").Value = 2.609843922559946
$ws.Range("H17").Value = 2.042167750953297
$ws.Range("I17").Value = 1.879283087165476
$ws.Range("J17").Value = 0.03597896206967377
$ws.Range("L17").Value = 0.6605018491596724
$ws.Range("B18").Value = 2.797538370780899
$ws.Range("C18").Value = 0.275712907937617
$ws.Range("D18").Value = 0.4476636562187224
$ws.Range("E18").Value = 0.1183674250960465
$ws.Range("G18").Value = 2.605118298697647
$ws.Range("H18").Value = 2.041717369511332
$ws.Range("I18").Value = 1.879687215680576
$ws.Range("J18").Value = 0.03587422017841746
$ws.Range("L18").Value = 0.6578257624865103
$ws.Range("B19").Value = 2.787191331552037
$ws.Range("C19").Value = 0.2735833718832907
$ws.Range("D19").Value = 0.4475186350129547
$ws.Range("E19").Value = 0.118404195650438
$ws.Range("G19").Value = 2.60354609188451
$ws.Range("H19").Value = 2.041582683918904
$ws.Range("I19").Value = 1.87984169581614
$ws.Range("J19").Value = 0.0358386866603162
$ws.Range("L19").Value = 0.6569248701981536
$ws.Range("B20").Value = 2.833831995780542
$ws.Range("C20").Value = 0.2831703481856778
$ws.Range("D20").Value = 0.4481839334041524
$ws.Range("E20").Value = 0.1182410006994203
$ws.Range("G20").Value = 2.610731750242167
$ws.Range("H20").Value = 2.042259568656135
$ws.Range("I20").Value = 1.879216681073885
$ws.Range("J20").Value = 0.03599831436032908
$ws.Range("L20").Value = 0.6609995956912087
$ws.Range("B21").Value = 2.991898784230216
$ws.Range("C21").Value = 0.3154489180666076
$ws.Range("D21").Value = 0.4506413101136957
$ws.Range("E21").Value = 0.1177324943705074
$ws.Range("G21").Value = 2.636816288816931
$ws.Range("H21").Value = 2.045775158264661
$ws.Range("I21").Value = 1.87834791039657
$ws.Range("J21").Value = 0.03652981175984493
$ws.Range("L21").Value = 0.6750509594356799
$ws.Range("B22").Value = 3.096216800478942
$ws.Range("C22").Value = 0.3365960763740077
$ws.Range("D22").Value = 0.4524129154284111
$ws.Range("E22").Value = 0.1174297824485429
$ws.Range("G22").Value = 2.655319781484906
$ws.Range("H22").Value = 2.049001183040957
$ws.Range("I22").Value = 1.878701801214447
$ws.Range("J22").Value = 0.03687359988033378
$ws.Range("L22").Value = 0.6845034994222772
$ws.Range("B23").Value = 3.040447800920674
$ws.Range("C23").Value = 0.3253046962329336
$ws.Range("D23").Value = 0.4514522227858606
$ws.Range("E23").Value = 0.1175886373452926
$ws.Range("G23").Value = 2.64531036725657
$ws.Range("H23").Value = 2.047194239653209
$ws.Range("I23").Value = 1.878428327349852
$ws.Range("J23").Value = 0.03669043937960481
$ws.Range("L23").Value = 0.6794338576990526
$ws.Range("B24").Value = 2.831265590025453
$ws.Range("C24").Value = 0.2826436248531934
$ws.Range("D24").Value = 0.448146563118641
$ws.Range("E24").Value = 0.1182498127258835
$ws.Range("G24").Value = 2.610329865111311
$ws.Range("H24").Value = 2.042217735421531
$ws.Range("I24").Value = 1.87924638237466
$ws.Range("J24").Value = 0.03598956659908303
$ws.Range("L24").Value = 0.6607744742530315
$ws.Range("B25").Value = 2.609985911450849
$ws.Range("C25").Value = 0.2368412985750297
$ws.Range("D25").Value = 0.4452913293299048
$ws.Range("E25").Value = 0.1190904683949903
$ws.Range("G25").Value = 2.578792221166594
$ws.Range("H25").Value = 2.040815002779482
$ws.Range("I25").Value = 1.884059625832521
$ws.Range("J25").Value = 0.03521805390600719
$ws.Range("L25").Value = 0.6418014661604303

Write-Output "Applied 216 cell updates"